$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    # Prefix with an apostrophe so Excel always stores the value as literal
    # text (never auto-converted to a number/date), then reset the style
    # back to Normal so no stray "quote prefix" cell format is left behind.
    $range.Value = "'" + $value
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "27.905.48"
Set-TextValue $ws.Range("E2") "  +1.28%  "

# Row 3
Set-TextValue $ws.Range("D3") "1.747.28"
Set-TextValue $ws.Range("E3") "  -0.87%  "

# Row 4
Set-TextValue $ws.Range("D4") "1.006"
Set-TextValue $ws.Range("E4") "  +0.09%  "

# Row 5
Set-TextValue $ws.Range("D5") "335.16"
Set-TextValue $ws.Range("E5") "  -0.24%  "

# Row 6
Set-TextValue $ws.Range("D6") "1.000"
Set-TextValue $ws.Range("E6") "  -0.07%  "

# Row 7
Set-TextValue $ws.Range("D7") "0.3844"
Set-TextValue $ws.Range("E7") "  +0.19%  "

# Row 8
Set-TextValue $ws.Range("D8") "0.3394"
Set-TextValue $ws.Range("E8") "  -0.40%  "

# Row 9
Set-TextValue $ws.Range("D9") "45.79"
Set-TextValue $ws.Range("E9") "  -2.64%  "

# Row 10
Set-TextValue $ws.Range("D10") "1.114"
Set-TextValue $ws.Range("E10") "  -2.21%  "

# Row 11
Set-TextValue $ws.Range("D11") "0.07197"
Set-TextValue $ws.Range("E11") "  -2.65%  "

# Row 12
Set-TextValue $ws.Range("B12") "BinanceUSD"
Set-TextValue $ws.Range("C12") "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextValue $ws.Range("D12") "1.004"
Set-TextValue $ws.Range("E12") "  +0.23%  "

# Row 13
Set-TextValue $ws.Range("B13") "Solana"
Set-TextValue $ws.Range("C13") "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-TextValue $ws.Range("D13") "22.47"
Set-TextValue $ws.Range("E13") "  -0.02%  "

# Row 14
Set-TextValue $ws.Range("D14") "6.149"
Set-TextValue $ws.Range("E14") "  -3.08%  "

# Row 15
Set-TextValue $ws.Range("B15") "WrappedEther"
Set-TextValue $ws.Range("C15") "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws.Range("D15") "1.751.61"
Set-TextValue $ws.Range("E15") "  -0.71%  "

# Row 16
Set-TextValue $ws.Range("B16") "Chainlink"
Set-TextValue $ws.Range("C16") "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue $ws.Range("D16") "7.094"
Set-TextValue $ws.Range("E16") "  +0.77%  "

# Row 17
Set-TextValue $ws.Range("D17") "0.00001056"
Set-TextValue $ws.Range("E17") "  -1.61%  "

# Row 18
Set-TextValue $ws.Range("D18") "0.06599"
Set-TextValue $ws.Range("E18") "  -0.94%  "

# Row 19
Set-TextValue $ws.Range("D19") "79.11"
Set-TextValue $ws.Range("E19") "  -3.76%  "

# Row 20
Set-TextValue $ws.Range("D20") "1.001"
Set-TextValue $ws.Range("E20") "  -0.08%  "

# Row 21
Set-TextValue $ws.Range("D21") "16.71"
Set-TextValue $ws.Range("E21") "  -3.64%  "

# Row 22
Set-TextValue $ws.Range("D22") "6.163"
Set-TextValue $ws.Range("E22") "  -3.29%  "

# Row 23
Set-TextValue $ws.Range("D23") "27.945.91"
Set-TextValue $ws.Range("E23") "  +1.40%  "

# Row 24
Set-TextValue $ws.Range("D24") "11.62"
Set-TextValue $ws.Range("E24") "  -3.35%  "

# Row 25
Set-TextValue $ws.Range("D25") "2.397"
Set-TextValue $ws.Range("E25") "  +0.47%  "

# Row 26
Set-TextValue $ws.Range("D26") "153.35"
Set-TextValue $ws.Range("E26") "  +0.55%  "

# Row 27
Set-TextValue $ws.Range("D27") "19.80"
Set-TextValue $ws.Range("E27") "  -3.81%  "

# Row 28
Set-TextValue $ws.Range("D28") "2.290"
Set-TextValue $ws.Range("E28") "  -5.08%  "

# Row 29
Set-TextValue $ws.Range("D29") "1.951.95"
Set-TextValue $ws.Range("E29") "  -0.59%  "

# Row 30
Set-TextValue $ws.Range("D30") "1.271"
Set-TextValue $ws.Range("E30") "  -10.59%  "

# Row 31
Set-TextValue $ws.Range("D31") "130.89"
Set-TextValue $ws.Range("E31") "  -2.72%  "

# Row 32
Set-TextValue $ws.Range("D32") "4.029"
Set-TextValue $ws.Range("E32") "  +1.90%  "

# Row 33
Set-TextValue $ws.Range("D33") "5.804"
Set-TextValue $ws.Range("E33") "  -4.92%  "

# Row 34
Set-TextValue $ws.Range("D34") "0.08808"
Set-TextValue $ws.Range("E34") "  +0.16%  "

# Row 35
Set-TextValue $ws.Range("D35") "12.13"
Set-TextValue $ws.Range("E35") "  -4.49%  "

# Row 36
Set-TextValue $ws.Range("D36") "1.539"
Set-TextValue $ws.Range("E36") "  +2.14%  "

# Row 37
Set-TextValue $ws.Range("D37") "0.6536"
Set-TextValue $ws.Range("E37") "  -3.45%  "

# Row 38
Set-TextValue $ws.Range("D38") "0.02275"
Set-TextValue $ws.Range("E38") "  -5.83%  "

# Row 39
Set-TextValue $ws.Range("D39") "5.109"
Set-TextValue $ws.Range("E39") "  -3.97%  "

# Row 40
Set-TextValue $ws.Range("D40") "0.06119"
Set-TextValue $ws.Range("E40") "  -3.05%  "

# Row 41
Set-TextValue $ws.Range("D41") "0.2089"
Set-TextValue $ws.Range("E41") "  -4.03%  "

# Row 42
Set-TextValue $ws.Range("D42") "1.206"
Set-TextValue $ws.Range("E42") "  -3.10%  "

# Row 43
Set-TextValue $ws.Range("D43") "7.983"
Set-TextValue $ws.Range("E43") "  -3.08%  "

# Row 44
Set-TextValue $ws.Range("D44") "1.000"
Set-TextValue $ws.Range("E44") "  -0.04%  "

# Row 45
Set-TextValue $ws.Range("D45") "13.73"
Set-TextValue $ws.Range("E45") "  -3.49%  "

# Row 46
Set-TextValue $ws.Range("D46") "3.834"
Set-TextValue $ws.Range("E46") "  +0.21%  "

# Row 47
Set-TextValue $ws.Range("D47") "0.6010"
Set-TextValue $ws.Range("E47") "  -3.73%  "

# Row 48
Set-TextValue $ws.Range("D48") "126.28"
Set-TextValue $ws.Range("E48") "  -3.42%  "

# Row 49
Set-TextValue $ws.Range("D49") "1.996"
Set-TextValue $ws.Range("E49") "  -3.90%  "

# Row 50
Set-TextValue $ws.Range("D50") "1.166"
Set-TextValue $ws.Range("E50") "  +1.78%  "

# Row 51
Set-TextValue $ws.Range("D51") "1.108"
Set-TextValue $ws.Range("E51") "  +4.37%  "
